$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'31.024.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Formula = "'1.966.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Formula = "'248.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Formula = "'0.4889"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Formula = "'0.2960"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").Formula = "'0.06835"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Formula = "'19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").Formula = "'1.970.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Formula = "'0.07778"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Formula = "'5.455"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Formula = "'0.7038"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.71%  "
$ws.Range("D16").Formula = "'287.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Formula = "'31.032.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D19").Formula = "'0.000007749"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Formula = "'2.222.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").Formula = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Formula = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Formula = "'6.622"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Formula = "'10.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.05%  "
$ws.Range("D26").Formula = "'170.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").Formula = "'20.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Formula = "'2.198"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").Formula = "'0.1071"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Formula = "'1.453"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Formula = "'4.828"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +18.80%  "
$ws.Range("D32").Formula = "'4.515"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("D33").Formula = "'0.05106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").Formula = "'0.7757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.04%  "
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Formula = "'2.734"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").Formula = "'0.02045"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").Formula = "'6.514"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.50%  "
$ws.Range("D40").Formula = "'2.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Formula = "'0.8898"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Formula = "'0.4498"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Formula = "'110.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Formula = "'73.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Formula = "'7.551"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.44%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Formula = "'9.518"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Formula = "'972.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.23%  "
$ws.Range("D49").Formula = "'0.1269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("D50").Formula = "'36.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").Formula = "'0.4106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
